$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "27.262.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -4.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.856.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -5.64%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -1.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "321.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  -0.37%  "

# Row 6
$ws.Range("E6").Value2 = "  -0.93%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4494"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -6.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3850"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  -5.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "47.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -11.36%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07880"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -7.35%  "

# Row 11
$ws.Range("E11").Value2 = "  -4.34%  "

# Row 12
$ws.Range("E12").Value2 = "  -5.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "1.859.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -7.17%  "

# Row 14
$ws.Range("B14").Value2 = "Polkadot"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.870"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -5.45%  "

# Row 15
$ws.Range("B15").Value2 = "Chainlink"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.148"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -6.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -1.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.00001030"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -3.87%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "85.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -6.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06524"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -1.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "16.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -9.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -1.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.495"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -6.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "27.268.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -4.59%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "10.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -6.72%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.263"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -1.58%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "2.077.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -6.95%  "

# Row 27
$ws.Range("E27").Value2 = "  -2.86%  "

# Row 28
$ws.Range("E28").Value2 = "  -3.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -5.83%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "5.429"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -7.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "120.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -3.76%  "

# Row 32
$ws.Range("B32").Value2 = "Stellar"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.09268"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -4.18%  "

# Row 33
$ws.Range("B33").Value2 = "ImmutableX"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.9356"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -5.30%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.466"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +0.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.566"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -3.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "5.282"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -6.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.02220"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -5.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.05978"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -4.29%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "1.207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  -4.02%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "8.255"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -10.04%  "

# Row 41
$ws.Range("E41").Value2 = "  -1.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.5898"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -5.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1878"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -2.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "10.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  -10.26%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  -6.90%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.5612"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -5.83%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "11.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -9.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "3.352"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -1.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "1.914"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  -7.71%  "

# Row 50
$ws.Range("E50").Value2 = "  -0.25%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "108.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -3.12%  "
